$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (date 39400 / year 2007), shifting all subsequent
# rows up by one, matching the upstream re-run of the naive forecaster.
$ws.Rows(2).Delete()

# Fix up the values that changed as a result of the forecaster re-run.
# Column C (y_0_forecast) value for row 3 was recomputed.
$ws.Range("C3").Value = -1.611885206309638

# Column E (y_1_forecast) values: the first four rows no longer have a
# forecast value, and the remaining rows received freshly recomputed values.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()

$ws.Range("E6").Value = 0.4163953164477929
$ws.Range("E7").Value = 1.653207170606596
$ws.Range("E8").Value = 1.270027657109818
$ws.Range("E9").Value = 1.579162878845075
$ws.Range("E10").Value = 1.452243308058287
$ws.Range("E11").Value = 2.068578555939404
$ws.Range("E12").Value = 1.651658474923545
$ws.Range("E13").Value = 0.115841687688345
$ws.Range("E14").Value = -2.092304328310923
$ws.Range("E15").Value = 1.533339625605379
$ws.Range("E16").Value = 0.492911192428136
$ws.Range("E17").Value = 0.2100922168233987
$ws.Range("E18").Value = 0.5208382580577098
